# Apply the edit described by the diff:
#  - column B header label changes from "类型" to "标题"
#  - column B data rows change from "TEXT" to "标题1" / "标题2"
#    (this introduces one brand-new shared string, "标题2")
#  - the active selection moves from C8 to B4
#  - the workbook window is resized smaller (windowWidth/windowHeight)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "标题"
$ws.Range("B2").Value = "标题1"
$ws.Range("B3").Value = "标题2"

# Move / collapse the active selection onto B4 (was C8)
$ws.Range("B4").Select()

# Shrink the saved window size (OOXML bookViews windowWidth/windowHeight are
# stored in twips = points * 20; ActiveWindow.Width/Height are in points).
$win = $excel.ActiveWindow
$win.Width = 27920 / 20
$win.Height = 11160 / 20

$wb.Save()
